# Update the cross labels in column A to distinguish replicate F1/F2 rows
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "F1a"
$ws.Range("A4").Value = "F1b"
$ws.Range("A5").Value = "F2a"
$ws.Range("A6").Value = "F2b"

# Move the active selection to G8, matching the author's final cursor position
$ws.Range("G8").Select()
